$d = $word.ActiveDocument

$endRange = $d.Content
$endRange.Collapse(0)

$xmlFragment = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body>
<w:p>
  <w:r><w:t xml:space="preserve">Resaltar un punto concreto en </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>scatter</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>plot</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve">: </w:t></w:r>
  <w:hyperlink r:id="rEditHyperlink1" w:history="1">
    <w:r>
      <w:rPr><w:rStyle w:val="Hipervnculo"/></w:rPr>
      <w:t>https://stackoverflow.com/questions/</w:t>
    </w:r>
    <w:r>
      <w:rPr><w:rStyle w:val="Hipervnculo"/></w:rPr>
      <w:t>3</w:t>
    </w:r>
    <w:r>
      <w:rPr><w:rStyle w:val="Hipervnculo"/></w:rPr>
      <w:t>8512485/highlight-specific-points-in-matplotlib-scatterplot</w:t>
    </w:r>
  </w:hyperlink>
</w:p>
<w:p/>
<w:p>
  <w:r><w:t xml:space="preserve">Usar </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>dataframes</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> diferentes en un mismo </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>plot</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve">: </w:t></w:r>
  <w:hyperlink r:id="rEditHyperlink2" w:history="1">
    <w:r>
      <w:rPr><w:rStyle w:val="Hipervnculo"/></w:rPr>
      <w:t>https://stackoverflow.com/questions/59766933/plots-different-columns-of-different-dataframe-in-one-plot-as-scatter-plot</w:t>
    </w:r>
  </w:hyperlink>
</w:p>
<w:p/>
<w:p>
  <w:r><w:t xml:space="preserve">Dibujar una línea dentro de </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>scatter</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve">: </w:t></w:r>
  <w:hyperlink r:id="rEditHyperlink3" w:history="1">
    <w:r>
      <w:rPr><w:rStyle w:val="Hipervnculo"/></w:rPr>
      <w:t>https://stackoverflow.com/questions/12981696/how-to-draw-line-inside-a-scatter-plot</w:t>
    </w:r>
  </w:hyperlink>
</w:p>
<w:p/>
<w:p>
  <w:pPr>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
    <w:t>Set_context</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
    <w:t xml:space="preserve"> de seaborn: </w:t>
  </w:r>
  <w:hyperlink r:id="rEditHyperlink4" w:history="1">
    <w:r>
      <w:rPr>
        <w:rStyle w:val="Hipervnculo"/>
        <w:lang w:val="en-US"/>
      </w:rPr>
      <w:t>https://seaborn.pydata.org/generated/seaborn.set_context.html</w:t>
    </w:r>
  </w:hyperlink>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="256">
<pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
<Relationship Id="rEditHyperlink1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://stackoverflow.com/questions/38512485/highlight-specific-points-in-matplotlib-scatterplot" TargetMode="External"/>
<Relationship Id="rEditHyperlink2" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://stackoverflow.com/questions/59766933/plots-different-columns-of-different-dataframe-in-one-plot-as-scatter-plot" TargetMode="External"/>
<Relationship Id="rEditHyperlink3" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://stackoverflow.com/questions/12981696/how-to-draw-line-inside-a-scatter-plot" TargetMode="External"/>
<Relationship Id="rEditHyperlink4" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://seaborn.pydata.org/generated/seaborn.set_context.html" TargetMode="External"/>
</Relationships>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$endRange.InsertXML($xmlFragment)

Write-Host "Webgrafia entries inserted"
